$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 2
$ws.Range("A2").Value = "Blocs béton manufacturés"
$ws.Range("B2").Value = 0.2
$ws.Range("C2").Value = "Maçonnerie - Enduit"
$ws.Range("D2").Value = 0.02

# Add new row 3
$ws.Range("A3").Value = "Maçonnerie - Béton"
$ws.Range("B3").Value = 0.2
$ws.Range("C3").Value = "Isolant Laine de bois"
$ws.Range("D3").Value = 0.08
